# Weekly update: insert 3 new daily price rows (Papa, Vega Modelo de Temuco)
# right after the existing row 336, pushing the old rows 337-431 down to 340-434.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 337 (shifts old rows 337..431 down to 340..434).
$ws.Rows("337:339").Insert()

# Common columns shared by every data row in this block.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100114001
$categoria = "Papa"
$kgOUnid   = 25
$clasif    = "Hortaliza"

# New row 337
$r = 337
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = [DateTime]"2021-09-24"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Asterix"
$ws.Cells.Item($r, 9).Value = "1a (guarda)"
$ws.Cells.Item($r, 10).Value = 500
$ws.Cells.Item($r, 11).Value = 9000
$ws.Cells.Item($r, 12).Value = 9000
$ws.Cells.Item($r, 13).Value = 9000
$ws.Cells.Item($r, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($r, 16).Value = 360
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif

# New row 338
$r = 338
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = [DateTime]"2021-09-24"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Patagonia"
$ws.Cells.Item($r, 9).Value = "1a (guarda)"
$ws.Cells.Item($r, 10).Value = 500
$ws.Cells.Item($r, 11).Value = 8000
$ws.Cells.Item($r, 12).Value = 8000
$ws.Cells.Item($r, 13).Value = 8000
$ws.Cells.Item($r, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($r, 16).Value = 320
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif

# New row 339
$r = 339
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = [DateTime]"2021-09-24"
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $catId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = "Rosara"
$ws.Cells.Item($r, 9).Value = "1a (guarda)"
$ws.Cells.Item($r, 10).Value = 600
$ws.Cells.Item($r, 11).Value = 8000
$ws.Cells.Item($r, 12).Value = 8000
$ws.Cells.Item($r, 13).Value = 8000
$ws.Cells.Item($r, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($r, 16).Value = 320
$ws.Cells.Item($r, 17).Value = $kgOUnid
$ws.Cells.Item($r, 18).Value = $clasif
